# "added Logic in Excel"
#
# 1. Delete the near-empty "SUMIFS" sheet (duplicate / leftover sheet).
# 2. On the "IFS" sheet, replace the boolean CRD flag column (C4:C13) with
#    text values ("CRD" / "No CRD" / "Nothing") that the IFS() formula in
#    column E actually compares against (C = "CRD"), and re-select the
#    edited range.
# 3. Nudge the selection / scroll position on a couple of other sheets to
#    match where the author ended up looking.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the stray "SUMIFS" sheet -----------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets("SUMIFS").Delete() | Out-Null
$excel.DisplayAlerts = $true

# --- 2. Fix up the CRD column on the "IFS" sheet ---------------------------
$ifs = $wb.Worksheets("IFS")

$ifs.Range("C4").Value = "CRD"
$ifs.Range("C5").Value = "CRD"
$ifs.Range("C6").Value = "No CRD"
$ifs.Range("C7").Value = "No CRD"
$ifs.Range("C8").Value = "CRD"
$ifs.Range("C9").Value = "No CRD"
$ifs.Range("C10").Value = "No CRD"
$ifs.Range("C11").Value = "No CRD"
$ifs.Range("C12").Value = "CRD"
$ifs.Range("C13").Value = "Nothing"

$ifs.Activate() | Out-Null
$ifs.Range("C4:C13").Select() | Out-Null

# --- 3. Match the author's final selections on a couple of sheets ---------
$sumifs = $wb.Worksheets("SUMIF, COUNTIF, SUMIFS")
$sumifs.Activate() | Out-Null
$sumifs.Range("I19").Select() | Out-Null

$examples = $wb.Worksheets("Example applications")
$examples.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$examples.Range("E45").Select() | Out-Null

$wb.Windows.Item(1).Left = 29850
$wb.Windows.Item(1).Top = -120
$wb.Windows.Item(1).Width = 27870
$wb.Windows.Item(1).Height = 16440
